$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

$ws.Range("N2:N12").Value = "nan"
